# Added SRRIP simulations for gromacs benchmark
# (also fills in the remaining gobmk/Hawkeye/OPTGen and gromacs/LRU/SRRIP
# rows that were previously blank and evaluating to #DIV/0!)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Config1 sheet - rows 24-28 (gobmk SRRIP/Hawkeye/OPTGen, gromacs LRU/SRRIP)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Config1")

# Row 24 - gobmk / SRRIP
$ws1.Cells.Item(24, 3).Value = 50000000
$ws1.Cells.Item(24, 4).Value = 116980411
$ws1.Cells.Item(24, 5).Value = 31351
$ws1.Cells.Item(24, 6).Value = 19240
$ws1.Cells.Item(24, 7).Value = 12111

# Row 25 - gobmk / Hawkeye
$ws1.Cells.Item(25, 3).Value = 50000000
$ws1.Cells.Item(25, 4).Value = 117019326
$ws1.Cells.Item(25, 5).Value = 31351
$ws1.Cells.Item(25, 6).Value = 17616
$ws1.Cells.Item(25, 7).Value = 13735

# Row 26 - gobmk / OPTGen
$ws1.Cells.Item(26, 3).Value = 50000000
$ws1.Cells.Item(26, 4).Value = 117019326
$ws1.Cells.Item(26, 5).Value = 1411
$ws1.Cells.Item(26, 6).Value = 859
$ws1.Cells.Item(26, 7).Formula = "=E26-F26"
$ws1.Cells.Item(26, 10).Formula = "=F26/E26"

# Row 27 - gromacs / LRU
$ws1.Cells.Item(27, 3).Value = 50000001
$ws1.Cells.Item(27, 4).Value = 124667669
$ws1.Cells.Item(27, 5).Value = 30089
$ws1.Cells.Item(27, 6).Value = 15882
$ws1.Cells.Item(27, 7).Value = 14207

# Row 28 - gromacs / SRRIP
$ws1.Cells.Item(28, 3).Value = 50000001
$ws1.Cells.Item(28, 4).Value = 124667669
$ws1.Cells.Item(28, 5).Value = 30089
$ws1.Cells.Item(28, 6).Value = 15632
$ws1.Cells.Item(28, 7).Value = 14457

# ---------------------------------------------------------------------
# Config2 sheet - rows 24-28 (same benchmark/policy layout as Config1)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Config2")

# Row 24 - gobmk / SRRIP
$ws2.Cells.Item(24, 3).Value = 50000000
$ws2.Cells.Item(24, 4).Value = 116540663
$ws2.Cells.Item(24, 5).Value = 46889
$ws2.Cells.Item(24, 6).Value = 28912
$ws2.Cells.Item(24, 7).Value = 17977

# Row 25 - gobmk / Hawkeye
$ws2.Cells.Item(25, 3).Value = 50000000
$ws2.Cells.Item(25, 4).Value = 116557215
$ws2.Cells.Item(25, 5).Value = 46888
$ws2.Cells.Item(25, 6).Value = 25562
$ws2.Cells.Item(25, 7).Value = 21326

# Row 26 - gobmk / OPTGen
$ws2.Cells.Item(26, 3).Value = 50000000
$ws2.Cells.Item(26, 4).Value = 116557215
$ws2.Cells.Item(26, 5).Value = 1167
$ws2.Cells.Item(26, 6).Value = 1035
$ws2.Cells.Item(26, 7).Formula = "=E26-F26"
$ws2.Cells.Item(26, 10).Formula = "=F26/E26"

# Row 27 - gromacs / LRU
$ws2.Cells.Item(27, 3).Value = 50000001
$ws2.Cells.Item(27, 4).Value = 124666538
$ws2.Cells.Item(27, 5).Value = 30134
$ws2.Cells.Item(27, 6).Value = 15919
$ws2.Cells.Item(27, 7).Value = 14215

# Row 28 - gromacs / SRRIP
$ws2.Cells.Item(28, 3).Value = 50000001
$ws2.Cells.Item(28, 4).Value = 124666582
$ws2.Cells.Item(28, 5).Value = 30134
$ws2.Cells.Item(28, 6).Value = 15661
$ws2.Cells.Item(28, 7).Value = 14473

# ---------------------------------------------------------------------
# View/selection state: Benchmarks keeps a new selection (C3), Config2
# loses the active tab, Config1 becomes the active tab with a new
# selection/scroll position (C29).
# ---------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Benchmarks")
$wsB.Range("C3").Select()

$ws2.Range("C29").Select()

$ws1.Range("C29").Select()
